# Insert a new data row before current row 360 (Apio / Macroferia Regional de Talca),
# shifting existing rows 360:389 down to 361:390, then populate the new row 360
# with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 360 downwards (and below) by inserting a new blank row at 360.
$ws.Rows.Item(360).Insert()

# Populate the newly inserted row 360 with the new record's data.
$ws.Cells.Item(360, 1).Value = 5
$ws.Cells.Item(360, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(360, 3).Value = "Maule"
$ws.Cells.Item(360, 4).Value = 45223
$ws.Cells.Item(360, 5).Value = 7
$ws.Cells.Item(360, 6).Value = 100112017
$ws.Cells.Item(360, 7).Value = "Apio"
$ws.Cells.Item(360, 8).Value = "Americana (o)"
$ws.Cells.Item(360, 9).Value = "Primera"
$ws.Cells.Item(360, 10).Value = 300
$ws.Cells.Item(360, 11).Value = 7500
$ws.Cells.Item(360, 12).Value = 7500
$ws.Cells.Item(360, 13).Value = 7500
$ws.Cells.Item(360, 14).Value = "`$/docena de matas"
$ws.Cells.Item(360, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(360, 16).Value = 1250
$ws.Cells.Item(360, 17).Value = 6
$ws.Cells.Item(360, 18).Value = "Hortaliza"

# Match the date-style formatting used by column D in the surrounding rows.
$ws.Cells.Item(360, 4).NumberFormat = $ws.Cells.Item(361, 4).NumberFormat
